$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 92, shifting existing rows 92:200 down to 93:201
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new data record
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100114014
$ws.Cells.Item(92, 7).Value = "Betarraga"
$ws.Cells.Item(92, 8).Value = "Sin especificar"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 500
$ws.Cells.Item(92, 11).Value = 1200
$ws.Cells.Item(92, 12).Value = 1200
$ws.Cells.Item(92, 13).Value = 1200
$ws.Cells.Item(92, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(92, 15).Value = "Región del Maule"
$ws.Cells.Item(92, 16).Value = 240
$ws.Cells.Item(92, 17).Value = 5
$ws.Cells.Item(92, 18).Value = "Hortaliza"
